$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (sharedStrings runs embedded in cells A8 and C9)
# ---------------------------------------------------------------------------
# A8: "Volume 31   Number  5" -> "Volume 31   Number  6"
$ws.Range("A8").Value = "Volume 31   Number  6"
# C9: "Report Covering the Week  1/29/2024  Through  2/4/2024"
#  -> "Report Covering the Week  2/5/2024  Through  2/11/2024"
$ws.Range("C9").Value = "Report Covering the Week  2/5/2024  Through  2/11/2024"

# ---------------------------------------------------------------------------
# Crime-stat table updates (rows 15-27)
# ---------------------------------------------------------------------------

# --- Row 15 (Rape) ---
# D15 and E15 flip from text ("0" / "***.*") to real numbers; pull number
# formatting from sibling cells that already carry the right style so we
# reuse existing style records instead of minting new ones.
$ws.Range("F15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1

$ws.Range("H15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100

$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 100
$ws.Range("M15").Value = 33.333333333333

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 16
$ws.Range("H16").Value = -15.789473684210
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 12.5
$ws.Range("M16").Value = -6.896551724137
$ws.Range("N16").Value = -53.448275862069

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 250
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 200
$ws.Range("I17").Value = 34
$ws.Range("J17").Value = 19
$ws.Range("K17").Value = 78.947368421052
$ws.Range("L17").Value = 25.925925925925
$ws.Range("M17").Value = 54.545454545454
$ws.Range("N17").Value = 9.677419354838

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 6
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 12
$ws.Range("J18").Value = 7
$ws.Range("K18").Value = 71.428571428571
$ws.Range("L18").Value = 33.333333333333
$ws.Range("M18").Value = -63.636363636363
$ws.Range("N18").Value = -88.785046728972

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 44.444444444444
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = 45
$ws.Range("I19").Value = 81
$ws.Range("J19").Value = 65
$ws.Range("K19").Value = 24.615384615384
$ws.Range("L19").Value = 37.288135593220
$ws.Range("M19").Value = 72.340425531914
$ws.Range("N19").Value = 37.288135593220

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 24
$ws.Range("E20").Value = -79.166666666666
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 49
$ws.Range("H20").Value = -48.979591836734
$ws.Range("I20").Value = 45
$ws.Range("J20").Value = 71
$ws.Range("K20").Value = -36.619718309859
$ws.Range("L20").Value = 25
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -79.545454545454

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -17.073170731707
$ws.Range("F21").Value = 131
$ws.Range("G21").Value = 124
$ws.Range("H21").Value = 5.645161290322
$ws.Range("I21").Value = 204
$ws.Range("J21").Value = 192
$ws.Range("K21").Value = 6.25
$ws.Range("L21").Value = 29.113924050632
$ws.Range("M21").Value = 36.912751677852
$ws.Range("N21").Value = -57.411273486430

# --- Row 23 (Housing) ---
# F23 flips from a number (2) to text "0" - force text formatting via the
# "@" number format, then restore the original General-text style (as used
# by the row's other text cells, e.g. C23) via a formats-only paste so we
# don't leave F23 on a stray numeric style.
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("F23").PasteSpecial(-4122)

$ws.Range("H23").Value = -100
$ws.Range("I23").Value = 6
$ws.Range("K23").Value = 200
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 50

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 23
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 94
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = 14.634146341463
$ws.Range("I24").Value = 131
$ws.Range("J24").Value = 137
$ws.Range("K24").Value = -4.379562043795
$ws.Range("L24").Value = 7.377049180327
$ws.Range("M24").Value = 3.149606299212

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 125
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -5.714285714285
$ws.Range("I25").Value = 49
$ws.Range("J25").Value = 49
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 2.083333333333
$ws.Range("M25").Value = 25.641025641025

# --- Row 26 (UCR Rape*) ---
$ws.Range("F16").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1

$ws.Range("H16").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = -100

$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = -42.857142857142

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 3
$ws.Range("I27").Value = 7
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = 40
$ws.Range("L27").Value = 75
